$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 6000
$ws.Range("I20").Value = 6000
$ws.Range("K20").Value = 6000
$ws.Range("M20").Value = -5770

$ws.Range("H35").Value = 6000
$ws.Range("I35").Value = 6000
$ws.Range("K35").Value = 6000
$ws.Range("M35").Value = -5621

$ws.Range("H70").Value = 3551.1614
$ws.Range("I70").Value = 1514.3334
$ws.Range("J70").Value = 4040
$ws.Range("K70").Value = 4543.0002
$ws.Range("L70").Value = 12120
$ws.Range("M70").Value = -4273.0002
$ws.Range("N70").Value = -12660

$ws.Range("H73").Value = 3551.1614
$ws.Range("I73").Value = 1514.3334
$ws.Range("J73").Value = 4040
$ws.Range("K73").Value = 4543.0002
$ws.Range("L73").Value = 12120
$ws.Range("M73").Value = -3607.0002
$ws.Range("N73").Value = -13992

$ws.Range("H98").Value = 930.61536
$ws.Range("I98").Value = 973
$ws.Range("J98").Value = 697.5
$ws.Range("K98").Value = 973
$ws.Range("L98").Value = 697.5
$ws.Range("M98").Value = 525
$ws.Range("N98").Value = -3693.5

$ws.Range("H122").Value = 930.61536
$ws.Range("I122").Value = 973
$ws.Range("J122").Value = 697.5
$ws.Range("K122").Value = 2919
$ws.Range("L122").Value = 2092.5
$ws.Range("M122").Value = -469
$ws.Range("N122").Value = -6992.5

$ws.Range("H127").Value = 1036.6666
$ws.Range("I127").Value = 1044
$ws.Range("K127").Value = 3132
$ws.Range("M127").Value = 1828

$ws.Range("H132").Value = 1583.3043
$ws.Range("I132").Value = 1353.1428
$ws.Range("K132").Value = 4059.4284
$ws.Range("M132").Value = -1529.4284

$ws.Range("H134").Value = 100250
$ws.Range("J134").Value = 100250
$ws.Range("L134").Value = 100250
$ws.Range("N134").Value = -110390

$ws.Range("H136").Value = 143842.4
$ws.Range("J136").Value = 143842.4
$ws.Range("L136").Value = 143842.4
$ws.Range("N136").Value = -154042.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1265.026
$ws.Range("I32").Value = 1268.5132
$ws.Range("K32").Value = 1268.5132
$ws.Range("M32").Value = -981.5132000000001

$ws.Range("H61").Value = 10158.4
$ws.Range("I61").Value = 9304
$ws.Range("J61").Value = 15000
$ws.Range("K61").Value = 9304
$ws.Range("L61").Value = 15000
$ws.Range("M61").Value = -9092
$ws.Range("N61").Value = -15424

$ws.Range("H74").Value = 2318.5667
$ws.Range("I74").Value = 1686.711
$ws.Range("J74").Value = 4214.1333
$ws.Range("K74").Value = 1686.711
$ws.Range("L74").Value = 4214.1333
$ws.Range("M74").Value = -812.711
$ws.Range("N74").Value = -5962.1333

$ws.Range("H77").Value = 2318.5667
$ws.Range("I77").Value = 1686.711
$ws.Range("J77").Value = 4214.1333
$ws.Range("K77").Value = 8433.555
$ws.Range("L77").Value = 21070.6665
$ws.Range("M77").Value = -4065.555
$ws.Range("N77").Value = -29806.6665

$ws.Range("H102").Value = 4348.68
$ws.Range("I102").Value = 3497
$ws.Range("K102").Value = 3497
$ws.Range("M102").Value = -1875

$ws.Range("H110").Value = 8360.308000000001
$ws.Range("I110").Value = 5460.625
$ws.Range("K110").Value = 5460.625
$ws.Range("M110").Value = -3415.625

$ws.Range("H136").Value = 10158.4
$ws.Range("I136").Value = 9304
$ws.Range("J136").Value = 15000
$ws.Range("K136").Value = 27912
$ws.Range("L136").Value = 45000
$ws.Range("M136").Value = -25362
$ws.Range("N136").Value = -50100

$ws.Range("H139").Value = 121061
$ws.Range("J139").Value = 121061
$ws.Range("L139").Value = 121061
$ws.Range("N139").Value = -131341

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3692.318
$ws.Range("I20").Value = 2903.7
$ws.Range("K20").Value = 2903.7
$ws.Range("M20").Value = -2656.7

$ws.Range("H99").Value = 4576.5386
$ws.Range("I99").Value = 3360.5557
$ws.Range("K99").Value = 3360.5557
$ws.Range("M99").Value = -1862.5557

$ws.Range("H134").Value = 5553.8667
$ws.Range("I134").Value = 5573.069
$ws.Range("J134").Value = 4997
$ws.Range("K134").Value = 16719.207
$ws.Range("L134").Value = 14991
$ws.Range("M134").Value = -14184.207
$ws.Range("N134").Value = -20061

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5497.8823
$ws.Range("I31").Value = 4230.933
$ws.Range("K31").Value = 4230.933
$ws.Range("M31").Value = -3935.933

$ws.Range("H34").Value = 5497.8823
$ws.Range("I34").Value = 4230.933
$ws.Range("K34").Value = 4230.933
$ws.Range("M34").Value = -4028.933

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").ClearContents()
$ws.Range("N119").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 8494.833000000001
$ws.Range("I23").Value = 199.14285
$ws.Range("K23").Value = 597.4285500000001
$ws.Range("M23").Value = -362.4285500000001

$ws.Range("H34").Value = 334.08334
$ws.Range("I34").Value = 170
$ws.Range("J34").Value = 498.16666
$ws.Range("K34").Value = 510
$ws.Range("L34").Value = 1494.49998
$ws.Range("M34").Value = -426
$ws.Range("N34").Value = -1662.49998

$ws.Range("H38").Value = 742.4286
$ws.Range("I38").Value = 46.666668
$ws.Range("J38").Value = 932.1818
$ws.Range("K38").Value = 140.000004
$ws.Range("L38").Value = 2796.5454
$ws.Range("M38").Value = 206.999996
$ws.Range("N38").Value = -3490.5454

$ws.Range("H39").Value = 1291.6666
$ws.Range("I39").Value = 900
$ws.Range("J39").Value = 1487.5
$ws.Range("K39").Value = 2700
$ws.Range("L39").Value = 4462.5
$ws.Range("M39").Value = -2406
$ws.Range("N39").Value = -5050.5

$ws.Range("H55").Value = 684.6
$ws.Range("I55").Value = 452
$ws.Range("J55").Value = 839.6667
$ws.Range("K55").Value = 1356
$ws.Range("L55").Value = 2519.0001
$ws.Range("M55").Value = -1179
$ws.Range("N55").Value = -2873.0001

$ws.Range("H86").Value = 675
$ws.Range("I86").Value = 550
$ws.Range("K86").Value = 1650
$ws.Range("M86").Value = -464

$ws.Range("H89").Value = 675
$ws.Range("I89").Value = 550
$ws.Range("K89").Value = 4950
$ws.Range("M89").Value = 978

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 657.6111
$ws.Range("I55").Value = 707.3333
$ws.Range("K55").Value = 707.3333
$ws.Range("M55").Value = -534.3333

$ws.Range("H122").Value = 4287.6665
$ws.Range("I122").Value = 4158.4
$ws.Range("K122").Value = 12475.2
$ws.Range("M122").Value = -10025.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 461.35
$ws.Range("I107").Value = 392.3125
$ws.Range("K107").Value = 1176.9375
$ws.Range("M107").Value = 743.0625

$ws.Range("H132").Value = 6498.05
$ws.Range("I132").Value = 5747.625
$ws.Range("K132").Value = 17242.875
$ws.Range("M132").Value = -14712.875
